$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename header row: "<Name>_old" -> "<Name>_FV2210", "<Name>_new" -> "<Name>_FV2304"
# ---------------------------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2210"
}
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2304"
}

# ---------------------------------------------------------------------
# 2. Freeze the header row (split below row 1)
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3. Turn the used range into an Excel Table ("Table1"), while preserving
#    the existing direct header-row formatting (bold/fill/border/align)
#    instead of letting the new table contribute its own header dxf.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$headerRange.Copy($scratch)
$headerRange.ClearFormats()

$range = $ws.Range("A1:U87")
$table = $ws.ListObjects.Add(1, $range, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$scratch.Clear()

$wb.Save()
